# Auto-generated Excel COM-interop script
# Refreshes market-price-derived columns (H-N) for specific leve rows
# across all 8 sheets, per the scheduled-runner data update.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").Value = $null
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").Value = $null
$ws.Range("H70").Value = 2313.2144
$ws.Range("I70").Value = 1841
$ws.Range("J70").Value = 2785.4285
$ws.Range("K70").Value = 5523
$ws.Range("L70").Value = 8356.2855
$ws.Range("M70").Value = -5253
$ws.Range("N70").Value = -8896.2855
$ws.Range("H73").Value = 2313.2144
$ws.Range("I73").Value = 1841
$ws.Range("J73").Value = 2785.4285
$ws.Range("K73").Value = 5523
$ws.Range("L73").Value = 8356.2855
$ws.Range("M73").Value = -4587
$ws.Range("N73").Value = -10228.2855
$ws.Range("H76").Value = 11885.5625
$ws.Range("I76").Value = 14216.9
$ws.Range("K76").Value = 14216.9
$ws.Range("M76").Value = -13901.9
$ws.Range("H79").Value = 11885.5625
$ws.Range("I79").Value = 14216.9
$ws.Range("K79").Value = 14216.9
$ws.Range("M79").Value = -13124.9
$ws.Range("H86").Value = 4875.222
$ws.Range("I86").Value = 4312.8335
$ws.Range("K86").Value = 4312.8335
$ws.Range("M86").Value = -3189.8335
$ws.Range("H89").Value = 4875.222
$ws.Range("I89").Value = 4312.8335
$ws.Range("K89").Value = 21564.1675
$ws.Range("M89").Value = -15948.1675
$ws.Range("H96").Value = 1882.625
$ws.Range("I96").Value = 1444.5294
$ws.Range("J96").Value = 2946.5715
$ws.Range("K96").Value = 4333.5882
$ws.Range("L96").Value = 8839.7145
$ws.Range("M96").Value = -2960.5882
$ws.Range("N96").Value = -11585.7145
$ws.Range("H98").Value = 1062
$ws.Range("I98").Value = 992.4706
$ws.Range("K98").Value = 992.4706
$ws.Range("M98").Value = 505.5294
$ws.Range("H107").Value = 2888
$ws.Range("I107").Value = 3400
$ws.Range("J107").Value = 2802.6667
$ws.Range("K107").Value = 3400
$ws.Range("L107").Value = 2802.6667
$ws.Range("M107").Value = -1480
$ws.Range("N107").Value = -6642.6667
$ws.Range("H122").Value = 1062
$ws.Range("I122").Value = 992.4706
$ws.Range("K122").Value = 2977.4118
$ws.Range("M122").Value = -527.4117999999999
$ws.Range("H125").Value = 5077.5557
$ws.Range("I125").Value = 4100
$ws.Range("J125").Value = 5356.857
$ws.Range("K125").Value = 36900
$ws.Range("L125").Value = 48211.713
$ws.Range("M125").Value = -34440
$ws.Range("N125").Value = -53131.713
$ws.Range("H135").Value = 14707219
$ws.Range("I135").Value = 931
$ws.Range("J135").Value = 41668748
$ws.Range("K135").Value = 8379
$ws.Range("L135").Value = 375018732
$ws.Range("M135").Value = -5844
$ws.Range("N135").Value = -375023802
$ws.Range("H137").Value = 2326.8108
$ws.Range("I137").Value = 2455.3462
$ws.Range("K137").Value = 7366.0386
$ws.Range("M137").Value = -4816.0386
$ws.Range("H138").Value = 10755024
$ws.Range("J138").Value = 19611092
$ws.Range("L138").Value = 58833276
$ws.Range("N138").Value = -58843556
$ws.Range("H141").Value = 2979.6365
$ws.Range("I141").Value = 2979.6365
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 8938.9095
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -3758.9095
$ws.Range("N141").Value = $null

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1103.25
$ws.Range("I2").Value = 877.8889
$ws.Range("J2").Value = 1393
$ws.Range("K2").Value = 877.8889
$ws.Range("L2").Value = 1393
$ws.Range("M2").Value = -764.8889
$ws.Range("N2").Value = -1619
$ws.Range("H32").Value = 11910820
$ws.Range("I32").Value = 15875343
$ws.Range("J32").Value = 17250.477
$ws.Range("K32").Value = 15875343
$ws.Range("L32").Value = 17250.477
$ws.Range("M32").Value = -15875056
$ws.Range("N32").Value = -17824.477
$ws.Range("H61").Value = 29414700
$ws.Range("I61").Value = 41667850
$ws.Range("K61").Value = 41667850
$ws.Range("M61").Value = -41667638
$ws.Range("H74").Value = 58890420
$ws.Range("I74").Value = 62570740
$ws.Range("J74").Value = 5300
$ws.Range("K74").Value = 62570740
$ws.Range("L74").Value = 5300
$ws.Range("M74").Value = -62569866
$ws.Range("N74").Value = -7048
$ws.Range("H77").Value = 58890420
$ws.Range("I77").Value = 62570740
$ws.Range("J77").Value = 5300
$ws.Range("K77").Value = 312853700
$ws.Range("L77").Value = 26500
$ws.Range("M77").Value = -312849332
$ws.Range("N77").Value = -35236
$ws.Range("H97").Value = 3830
$ws.Range("I97").Value = 2995
$ws.Range("K97").Value = 2995
$ws.Range("M97").Value = -2499
$ws.Range("H102").Value = 3479.75
$ws.Range("I102").Value = 2537.5625
$ws.Range("J102").Value = 7248.5
$ws.Range("K102").Value = 2537.5625
$ws.Range("L102").Value = 7248.5
$ws.Range("M102").Value = -915.5625
$ws.Range("N102").Value = -10492.5
$ws.Range("H116").Value = 1103.25
$ws.Range("I116").Value = 877.8889
$ws.Range("J116").Value = 1393
$ws.Range("K116").Value = 877.8889
$ws.Range("L116").Value = 1393
$ws.Range("M116").Value = 1416.1111
$ws.Range("N116").Value = -5981
$ws.Range("H132").Value = 18520008
$ws.Range("I132").Value = 1421.3269
$ws.Range("J132").Value = 500003260
$ws.Range("K132").Value = 4263.9807
$ws.Range("L132").Value = 1500009780
$ws.Range("M132").Value = -1733.9807
$ws.Range("N132").Value = -1500014840
$ws.Range("H136").Value = 29414700
$ws.Range("I136").Value = 41667850
$ws.Range("K136").Value = 125003550
$ws.Range("M136").Value = -125001000

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1103.25
$ws.Range("I3").Value = 877.8889
$ws.Range("J3").Value = 1393
$ws.Range("K3").Value = 877.8889
$ws.Range("L3").Value = 1393
$ws.Range("M3").Value = -763.8889
$ws.Range("N3").Value = -1621
$ws.Range("H20").Value = 2141.8948
$ws.Range("I20").Value = 2170.6155
$ws.Range("K20").Value = 2170.6155
$ws.Range("M20").Value = -1923.6155
$ws.Range("H64").Value = 1738.3889
$ws.Range("I64").Value = 1453.3636
$ws.Range("K64").Value = 1453.3636
$ws.Range("M64").Value = -1228.3636
$ws.Range("H67").Value = 1738.3889
$ws.Range("I67").Value = 1453.3636
$ws.Range("K67").Value = 1453.3636
$ws.Range("M67").Value = -673.3635999999999
$ws.Range("H107").Value = 3860.9
$ws.Range("I107").Value = 3926.125
$ws.Range("J107").Value = 3600
$ws.Range("K107").Value = 3926.125
$ws.Range("L107").Value = 3600
$ws.Range("M107").Value = -2006.125
$ws.Range("N107").Value = -7440
$ws.Range("H134").Value = 3978.7144
$ws.Range("I134").Value = 3900.1924
$ws.Range("K134").Value = 11700.5772
$ws.Range("M134").Value = -9165.5772

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 5453.2383
$ws.Range("J22").Value = 1451.3334
$ws.Range("L22").Value = 1451.3334
$ws.Range("N22").Value = -2151.3334
$ws.Range("H31").Value = 24395086
$ws.Range("I31").Value = 3265.853
$ws.Range("J31").Value = 142869650
$ws.Range("K31").Value = 3265.853
$ws.Range("L31").Value = 142869650
$ws.Range("M31").Value = -2970.853
$ws.Range("N31").Value = -142870240
$ws.Range("H34").Value = 24395086
$ws.Range("I34").Value = 3265.853
$ws.Range("J34").Value = 142869650
$ws.Range("K34").Value = 3265.853
$ws.Range("L34").Value = 142869650
$ws.Range("M34").Value = -3063.853
$ws.Range("N34").Value = -142870054
$ws.Range("H94").Value = 1983.7778
$ws.Range("I94").Value = 1900
$ws.Range("J94").Value = 2007.7142
$ws.Range("K94").Value = 1900
$ws.Range("L94").Value = 2007.7142
$ws.Range("M94").Value = -1449
$ws.Range("N94").Value = -2909.7142
$ws.Range("H122").Value = 921.4
$ws.Range("I122").Value = 929.73914
$ws.Range("J122").Value = 905.4167
$ws.Range("K122").Value = 2789.21742
$ws.Range("L122").Value = 2716.2501
$ws.Range("M122").Value = -339.2174199999999
$ws.Range("N122").Value = -7616.2501
$ws.Range("H132").Value = 3224.1072
$ws.Range("I132").Value = 2521
$ws.Range("J132").Value = 5802.1665
$ws.Range("K132").Value = 7563
$ws.Range("L132").Value = 17406.4995
$ws.Range("M132").Value = -5033
$ws.Range("N132").Value = -22466.4995
$ws.Range("H134").Value = 1161.8636
$ws.Range("I134").Value = 1161.8636
$ws.Range("K134").Value = 3485.5908
$ws.Range("M134").Value = -950.5907999999999
$ws.Range("H141").Value = 262454.56
$ws.Range("I141").Value = 48431.668
$ws.Range("J141").Value = 320824.47
$ws.Range("K141").Value = 48431.668
$ws.Range("L141").Value = 320824.47
$ws.Range("M141").Value = -43251.668
$ws.Range("N141").Value = -331184.47

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 50809044
$ws.Range("I4").Value = 92076370
$ws.Range("J4").Value = 22437762
$ws.Range("K4").Value = 276229110
$ws.Range("L4").Value = 67313286
$ws.Range("M4").Value = -276228998
$ws.Range("N4").Value = -67313510
$ws.Range("H62").Value = 4078
$ws.Range("I62").Value = 4078
$ws.Range("K62").Value = 12234
$ws.Range("M62").Value = -11548
$ws.Range("H65").Value = 4078
$ws.Range("I65").Value = 4078
$ws.Range("K65").Value = 36702
$ws.Range("M65").Value = -33270
$ws.Range("H86").Value = 677.36365
$ws.Range("I86").Value = 477.85715
$ws.Range("K86").Value = 1433.57145
$ws.Range("M86").Value = -247.5714499999999
$ws.Range("H89").Value = 677.36365
$ws.Range("I89").Value = 477.85715
$ws.Range("K89").Value = 4300.71435
$ws.Range("M89").Value = 1627.28565
$ws.Range("H122").Value = 1637.2307
$ws.Range("I122").Value = 1561.75
$ws.Range("K122").Value = 14055.75
$ws.Range("M122").Value = -11605.75
$ws.Range("H134").Value = 7561.28
$ws.Range("I134").Value = 1884.2354
$ws.Range("J134").Value = 19625
$ws.Range("K134").Value = 5652.706200000001
$ws.Range("L134").Value = 58875
$ws.Range("M134").Value = -582.7062000000005
$ws.Range("N134").Value = -69015
$ws.Range("H140").Value = 1736.125
$ws.Range("I140").Value = 1550.8572
$ws.Range("J140").Value = 3033
$ws.Range("K140").Value = 4652.571599999999
$ws.Range("L140").Value = 9099
$ws.Range("M140").Value = 527.4284000000007
$ws.Range("N140").Value = -19459

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 8250
$ws.Range("I43").Value = 8000
$ws.Range("J43").Value = 8500
$ws.Range("K43").Value = 8000
$ws.Range("L43").Value = 8500
$ws.Range("M43").Value = -7849
$ws.Range("N43").Value = -8802
$ws.Range("H46").Value = 19999.75
$ws.Range("I46").Value = 10000
$ws.Range("J46").Value = 23333
$ws.Range("K46").Value = 10000
$ws.Range("L46").Value = 23333
$ws.Range("M46").Value = -9844
$ws.Range("N46").Value = -23645
$ws.Range("H102").Value = 1668.3448
$ws.Range("I102").Value = 1210.7826
$ws.Range("J102").Value = 3422.3333
$ws.Range("K102").Value = 1210.7826
$ws.Range("L102").Value = 3422.3333
$ws.Range("M102").Value = 411.2174
$ws.Range("N102").Value = -6666.3333
$ws.Range("H132").Value = 2410.2307
$ws.Range("I132").Value = 2393.5
$ws.Range("K132").Value = 7180.5
$ws.Range("M132").Value = -4650.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2059.4546
$ws.Range("I22").Value = 1276.25
$ws.Range("J22").Value = 2712.125
$ws.Range("K22").Value = 1276.25
$ws.Range("L22").Value = 2712.125
$ws.Range("M22").Value = -981.25
$ws.Range("N22").Value = -3302.125
$ws.Range("H27").Value = 2059.4546
$ws.Range("I27").Value = 1276.25
$ws.Range("J27").Value = 2712.125
$ws.Range("K27").Value = 1276.25
$ws.Range("L27").Value = 2712.125
$ws.Range("M27").Value = -1169.25
$ws.Range("N27").Value = -2926.125
$ws.Range("H40").Value = 3398.739
$ws.Range("I40").Value = 2592.3333
$ws.Range("K40").Value = 2592.3333
$ws.Range("M40").Value = -2456.3333
$ws.Range("H55").Value = 577.8929000000001
$ws.Range("I55").Value = 398.26666
$ws.Range("J55").Value = 785.1539
$ws.Range("K55").Value = 398.26666
$ws.Range("L55").Value = 785.1539
$ws.Range("M55").Value = -225.26666
$ws.Range("N55").Value = -1131.1539
$ws.Range("H93").Value = 2099.25
$ws.Range("I93").Value = 999.1667
$ws.Range("K93").Value = 999.1667
$ws.Range("M93").Value = 248.8333
$ws.Range("H122").Value = 3913.476
$ws.Range("I122").Value = 3279.0667
$ws.Range("J122").Value = 5499.5
$ws.Range("K122").Value = 9837.2001
$ws.Range("L122").Value = 16498.5
$ws.Range("M122").Value = -7387.2001
$ws.Range("N122").Value = -21398.5
$ws.Range("H132").Value = 90912400
$ws.Range("I132").Value = 2999.8667
$ws.Range("K132").Value = 8999.6001
$ws.Range("M132").Value = -6469.6001
$ws.Range("H136").Value = 2193.606
$ws.Range("I136").Value = 1680
$ws.Range("J136").Value = 7329.6665
$ws.Range("K136").Value = 5040
$ws.Range("L136").Value = 21988.9995
$ws.Range("M136").Value = -2490
$ws.Range("N136").Value = -27088.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 112223230
$ws.Range("I100").Value = 144286450
$ws.Range("K100").Value = 288572900
$ws.Range("M100").Value = -288572359
$ws.Range("H113").Value = 668.5
$ws.Range("I113").Value = 379.875
$ws.Range("J113").Value = 957.125
$ws.Range("K113").Value = 1139.625
$ws.Range("L113").Value = 2871.375
$ws.Range("M113").Value = 1030.375
$ws.Range("N113").Value = -7211.375
$ws.Range("H122").Value = 47668204
$ws.Range("I122").Value = 52685436
$ws.Range("K122").Value = 158056308
$ws.Range("M122").Value = -158053858
$ws.Range("H126").Value = 3188.1562
$ws.Range("I126").Value = 3188.1562
$ws.Range("K126").Value = 9564.4686
$ws.Range("M126").Value = -7094.4686
$ws.Range("H132").Value = 5180.175
$ws.Range("I132").Value = 5171.8237
$ws.Range("J132").Value = 5227.5
$ws.Range("K132").Value = 15515.4711
$ws.Range("L132").Value = 15682.5
$ws.Range("M132").Value = -12985.4711
$ws.Range("N132").Value = -20742.5
$ws.Range("H136").Value = 1480
$ws.Range("I136").Value = 1099.3077
$ws.Range("J136").Value = 4779.3335
$ws.Range("K136").Value = 3297.9231
$ws.Range("L136").Value = 14338.0005
$ws.Range("M136").Value = -747.9231
$ws.Range("N136").Value = -19438.0005
